# Apply the "updated adj study calculated mean and sd" edit to the
# PTG effect-sizes workbook.
#
# Row 2  = Adjorlolo 2022: recalculated effect size (mean) and added sd;
#          comment resolved to "this is correct".
# Row 3  = Arnout 2021: fixed scale-type label (was "PTG-21", should be "PTGI").
# Row 8  = Feingold 2022: placeholder effect-size note cleared (blank pending data).
# Row 12 = Lewis 2021: added the missing sd value.
# Row 16 = Pirtrzak 2021: placeholder effect-size note cleared (blank pending data).
#
# Also updates the saved cursor/selection position to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Adjorlolo 2022: updated calculated mean and sd
$ws.Range("C2").Value = 22.258
$ws.Range("D2").Value = 5.052
$ws.Range("F2").Value = "this is correct"

# Row 3 - Arnout 2021: scale type correction
$ws.Range("B3").Value = "PTGI"

# Row 8 - Feingold 2022: clear placeholder effect-size text
$ws.Range("C8").ClearContents()

# Row 12 - Lewis 2021: add missing sd
$ws.Range("D12").Value = 11.01

# Row 16 - Pirtrzak 2021: clear placeholder effect-size text
$ws.Range("C16").ClearContents()

# Restore the cursor / active-cell selection
$ws.Range("B4").Select()
